# NYPD CompStat weekly report refresh: new crime data for the week of
# 6/3/2024 - 6/9/2024 (Volume 31, Number 23), replacing the prior week's
# figures (Volume 31, Number 22; week of 5/27/2024 - 6/2/2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text: issue number and report date range (rich-text cells).
# ---------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 31   Number  23"
$ws.Range("C9").Value  = "Report Covering the Week  6/3/2024  Through  6/9/2024"

# ---------------------------------------------------------------------
# 2. Bulk numeric refresh for the crime-category table (rows 15-31).
#    Most cells simply get a new number written over the old one - the
#    cell's existing numeric style/format is left untouched.
# ---------------------------------------------------------------------
$numericUpdates = @{
    "F15" = 3;               "G15" = 2;                "H15" = 50;                 "N15" = -78.787878787878;

    "C16" = 8;                "E16" = 166.666666666667; "F16" = 30;                 "G16" = 12;
    "H16" = 150;               "I16" = 107;              "J16" = 74;                 "K16" = 44.594594594594;
    "L16" = 42.666666666666;  "M16" = 3.883495145631;   "N16" = -68.343195266272;

    "C17" = 12;                "D17" = 16;               "E17" = -25;                "F17" = 45;
    "G17" = 38;                "H17" = 18.421052631578;  "I17" = 199;                "J17" = 165;
    "K17" = 20.606060606060;  "L17" = 31.788079470198;  "M17" = 158.441558441558;   "N17" = -18.775510204081;

    "C18" = 3;                 "F18" = 12;               "G18" = 11;                 "H18" = 9.090909090909;
    "I18" = 79;                "K18" = 51.923076923076;  "L18" = -1.25;              "M18" = 229.166666666667;
    "N18" = -53.801169590643;

    "F19" = 29;                "G19" = 33;               "H19" = -12.121212121212;  "I19" = 165;
    "J19" = 183;               "K19" = -9.836065573770;  "L19" = -16.243654822335;  "M19" = 58.653846153846;
    "N19" = -38.202247191011;

    "C20" = 2;                 "F20" = 6;                "G20" = 1;                  "H20" = 500;
    "I20" = 26;                "K20" = -31.578947368421; "L20" = -35;                "M20" = 44.444444444444;
    "N20" = -84.337349397590;

    "C21" = 31;                "D21" = 28;               "E21" = 10.714285714285;   "F21" = 126;
    "G21" = 97;                "H21" = 29.896907216494;  "I21" = 586;                "J21" = 526;
    "K21" = 11.406844106463;  "L21" = 5.395683453237;   "M21" = 75.449101796407;   "N21" = -52.550607287449;

    "I22" = 4;                 "J22" = 2;                "K22" = 100;                "L22" = 100;
    "M22" = 100;

    "D23" = 8;                 "E23" = 25;               "F23" = 44;                 "G23" = 33;
    "H23" = 33.333333333333;  "I23" = 175;               "J23" = 164;                "K23" = 6.707317073170;
    "L23" = 6.060606060606;   "M23" = 84.210526315789;

    "C24" = 9;                 "D24" = 16;               "E24" = -43.75;             "F24" = 48;
    "G24" = 58;                "H24" = -17.241379310344; "I24" = 333;                "J24" = 384;
    "K24" = -13.28125;         "L24" = -2.346041055718;  "M24" = 22.878228782287;

    "C25" = 4;                 "D25" = 7;                "E25" = -42.857142857142;  "F25" = 16;
    "G25" = 21;                "H25" = -23.809523809523; "I25" = 74;                 "J25" = 102;
    "K25" = -27.450980392156; "L25" = -13.953488372093;

    "C26" = 18;                "D26" = 18;               "E26" = 0;                  "F26" = 93;
    "G26" = 57;                "H26" = 63.157894736842;  "I26" = 351;                "J26" = 249;
    "K26" = 40.963855421686;  "L26" = 33.969465648855;  "M26" = 12.5;

    "D27" = 1;                 "E27" = -100;             "F27" = 4;                  "H27" = 0;
    "J27" = 16;                "K27" = -25;

    "C28" = 6;                 "D28" = 1;                "E28" = 500;                "F28" = 10;
    "G28" = 6;                 "H28" = 66.666666666666;  "I28" = 29;                 "J28" = 23;
    "K28" = 26.086956521739;  "L28" = 7.407407407407;

    "C29" = 1;                 "D29" = 1;                "E29" = 0;                  "F29" = 1;
    "G29" = 1;                 "H29" = 0;                "I29" = 4;                  "J29" = 9;
    "K29" = -55.555555555555; "L29" = -63.636363636363; "M29" = -73.333333333333;  "N29" = -89.189189189189;

    "C30" = 1;                 "D30" = 1;                "E30" = 0;                  "F30" = 1;
    "H30" = 0;                 "I30" = 3;                "J30" = 8;                  "K30" = -62.5;
    "L30" = -62.5;             "M30" = -76.923076923076; "N30" = -91.176470588235;

    "L31" = -100;
}

foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).Value = $numericUpdates[$addr]
}

# ---------------------------------------------------------------------
# 3. Cells that flip from a text placeholder ("0" / "***.*") to a real
#    number this week - write the value AND apply the matching numeric
#    format so the style index matches a genuine numeric cell.
# ---------------------------------------------------------------------
$toNumberPlain = @("C22", "D22", "F22", "G22", "C28", "C29", "D29", "F29", "C30", "D30", "F30")
foreach ($addr in $toNumberPlain) {
    $ws.Range($addr).NumberFormat = "#,##0"
}

$toNumberPct = @("E22", "H22", "E29", "E30")
foreach ($addr in $toNumberPct) {
    $ws.Range($addr).NumberFormat = "#,##0.0;""-""#,##0.0"
}

$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0

$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1

$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1

# ---------------------------------------------------------------------
# 4. Cells that flip the other way - from a real number to the text
#    placeholder "0" - this week's count dropped to zero/unreported.
#    Use a quote-prefixed value (forces text) then copy the number
#    format from an existing placeholder cell so the style matches.
# ---------------------------------------------------------------------
$zeroTemplate = $ws.Range("C14")   # existing text-style "0" placeholder cell

$ws.Range("D18").Value = "'0"
$zeroTemplate.Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$zeroTemplate.Copy()
$ws.Range("C27").PasteSpecial(-4122)

$starTemplate = $ws.Range("E14")   # existing text-style "***.*" placeholder cell

$ws.Range("E18").Value = "'***.*"
$starTemplate.Copy()
$ws.Range("E18").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5. Column H narrows this week to match the other percent-change
#    columns now that its values are shorter.
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth
